$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Auto TestPN 1"
$ws.Range("A3").Value = "Auto TestPN 2"
$ws.Range("A4").Value = "Auto TestPN 3"
